$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("VENTA MENSUAL")

# Update INODOROS/PRESUPUESTO (column G) values for rows with new sales figures,
# and refresh the PRESUPUESTO total in row 56.
$values = @{
    2 = 1500
    3 = 1500
    4 = 1000
    6 = 3000
    9 = 2000
    10 = 7000
    12 = 3000
    16 = 3000
    19 = 1000
    26 = 5000
    29 = 2000
    30 = 1000
    31 = 1000
    36 = 500
    37 = 500
    38 = 1000
    42 = 2000
    43 = 5000
    47 = 1000
    48 = 1000
    49 = 3000
    53 = 500
    54 = 3000
    55 = 1000
    56 = 50500
}

foreach ($row in $values.Keys) {
    $ws.Cells.Item($row, 7).Value = $values[$row]
}

